$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 44, which shifts the existing rows 44 and 45
# down to become rows 45 and 46 (with all their data/styles intact).
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record's data.
$ws.Range("A44").Value = 11
$ws.Range("B44").Value = "Vega Monumental Concepción"
$ws.Range("C44").Value = "Bíobío"
$ws.Range("D44").Value = 44706
$ws.Range("E44").Value = 8
$ws.Range("F44").Value = 100112030
$ws.Range("G44").Value = "Poroto granado"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 100
$ws.Range("K44").Value = 23000
$ws.Range("L44").Value = 25000
$ws.Range("M44").Value = 24000
$ws.Range("N44").Value = "$/saco 25 kilos"
$ws.Range("O44").Value = "Región Metropolitana"
$ws.Range("P44").Value = 960
$ws.Range("Q44").Value = 25
$ws.Range("R44").Value = "Hortaliza"
